# Add "hydrogen combined cycle" as a new power plant type row on the
# BZECfNP sheet, and rename the existing "hydrogen" row to
# "hydrogen combustion turbine".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BZECfNP")

# Rename existing "hydrogen" entry to "hydrogen combustion turbine"
$ws.Range("A24").Value = "hydrogen combustion turbine"

# Add the new "hydrogen combined cycle" row right below, mirroring the
# all-zero data of the row above (columns B:AE, years 2021-2047).
$ws.Range("A25").Value = "hydrogen combined cycle"
$ws.Range("B25:AE25").Value = 0

# Reflect the new row in the sheet's view (select it, scroll it into view),
# then restore "About" as the active sheet/tab, matching the saved file.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 4
$ws.Range("B25:AE25").Select()
$wb.Worksheets.Item("About").Activate()

